$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 142857820
$ws.Range("I32").Value = 333333800
$ws.Range("J32").Value = 845
$ws.Range("K32").Value = 333333800
$ws.Range("L32").Value = 845
$ws.Range("M32").Value = -333333474
$ws.Range("N32").Value = -1497

# Row 70
$ws.Range("H70").Value = 1692.3077
$ws.Range("I70").Value = 1000
$ws.Range("K70").Value = 3000
$ws.Range("M70").Value = -2730

# Row 73
$ws.Range("H73").Value = 1692.3077
$ws.Range("I73").Value = 1000
$ws.Range("K73").Value = 3000
$ws.Range("M73").Value = -2064

# Row 74
$ws.Range("H74").Value = 3797.1428
$ws.Range("I74").Value = 3485.3845
$ws.Range("J74").Value = 4303.75
$ws.Range("K74").Value = 3485.3845
$ws.Range("L74").Value = 4303.75
$ws.Range("M74").Value = -2549.3845
$ws.Range("N74").Value = -6175.75

# Row 77
$ws.Range("H77").Value = 3797.1428
$ws.Range("I77").Value = 3485.3845
$ws.Range("J77").Value = 4303.75
$ws.Range("K77").Value = 17426.9225
$ws.Range("L77").Value = 21518.75
$ws.Range("M77").Value = -12746.9225
$ws.Range("N77").Value = -30878.75

# Row 137
$ws.Range("H137").Value = 1103
$ws.Range("I137").Value = 962.4815
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 2887.4445
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -337.4445000000001
$ws.Range("N137").Value = -14100

# Row 138
$ws.Range("H138").Value = 2690.4375
$ws.Range("I138").Value = 1680.0416
$ws.Range("J138").Value = 3296.675
$ws.Range("K138").Value = 5040.1248
$ws.Range("L138").Value = 9890.025000000001
$ws.Range("M138").Value = 99.8752000000004
$ws.Range("N138").Value = -20170.025

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 28867.889
$ws.Range("I32").Value = 5084.0444
$ws.Range("J32").Value = 147787.11
$ws.Range("K32").Value = 5084.0444
$ws.Range("L32").Value = 147787.11
$ws.Range("M32").Value = -4797.0444
$ws.Range("N32").Value = -148361.11

# Row 61
$ws.Range("H61").Value = 4064.4
$ws.Range("I61").Value = 3965.2
$ws.Range("J61").Value = 4262.8
$ws.Range("K61").Value = 3965.2
$ws.Range("L61").Value = 4262.8
$ws.Range("M61").Value = -3753.2
$ws.Range("N61").Value = -4686.8

# Row 136
$ws.Range("H136").Value = 4064.4
$ws.Range("I136").Value = 3965.2
$ws.Range("J136").Value = 4262.8
$ws.Range("K136").Value = 11895.6
$ws.Range("L136").Value = 12788.4
$ws.Range("M136").Value = -9345.599999999999
$ws.Range("N136").Value = -17888.4

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1820.9
$ws.Range("I134").Value = 1519.65
$ws.Range("J134").Value = 2423.4
$ws.Range("K134").Value = 4558.950000000001
$ws.Range("L134").Value = 7270.200000000001
$ws.Range("M134").Value = -2023.950000000001
$ws.Range("N134").Value = -12340.2

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 17859052
$ws.Range("J31").Value = 2383.6453
$ws.Range("L31").Value = 2383.6453
$ws.Range("N31").Value = -2973.6453

# Row 34
$ws.Range("H34").Value = 17859052
$ws.Range("J34").Value = 2383.6453
$ws.Range("L34").Value = 2383.6453
$ws.Range("N34").Value = -2787.6453

# Row 43
$ws.Range("H43").Value = 29885.666
$ws.Range("J43").Value = 29885.666
$ws.Range("L43").Value = 29885.666
$ws.Range("N43").Value = -30253.666

# Row 99
$ws.Range("H99").Value = 1497.2142
$ws.Range("I99").Value = 1346.1
$ws.Range("J99").Value = 1875
$ws.Range("K99").Value = 1346.1
$ws.Range("L99").Value = 1875
$ws.Range("M99").Value = 151.9000000000001
$ws.Range("N99").Value = -4871

# Row 101
$ws.Range("H101").Value = 29885.666
$ws.Range("J101").Value = 29885.666
$ws.Range("L101").Value = 29885.666
$ws.Range("N101").Value = -36375.666

# Row 126
$ws.Range("H126").Value = 1497.2142
$ws.Range("I126").Value = 1346.1
$ws.Range("J126").Value = 1875
$ws.Range("K126").Value = 4038.3
$ws.Range("L126").Value = 5625
$ws.Range("M126").Value = -1568.3
$ws.Range("N126").Value = -10565

$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 1827.6471
$ws.Range("I46").Value = 287.5
$ws.Range("J46").Value = 3196.6667
$ws.Range("K46").Value = 862.5
$ws.Range("L46").Value = 9590.000100000001
$ws.Range("M46").Value = -771.5
$ws.Range("N46").Value = -9772.000100000001

# Row 68
$ws.Range("H68").Value = 1491.6578
$ws.Range("J68").Value = 1736.32
$ws.Range("L68").Value = 5208.96
$ws.Range("N68").Value = -6830.96

# Row 71
$ws.Range("H71").Value = 1491.6578
$ws.Range("J71").Value = 1736.32
$ws.Range("L71").Value = 15626.88
$ws.Range("N71").Value = -23738.88

# Row 107
$ws.Range("H107").Value = 642.1389
$ws.Range("I107").Value = 363.725
$ws.Range("J107").Value = 990.15625
$ws.Range("K107").Value = 1091.175
$ws.Range("L107").Value = 2970.46875
$ws.Range("M107").Value = 828.8249999999998
$ws.Range("N107").Value = -6810.46875

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 5820
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 8640
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 25920
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -30860

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1886.6666
$ws.Range("I68").Value = 1656
$ws.Range("J68").Value = 2175
$ws.Range("K68").Value = 1656
$ws.Range("L68").Value = 2175
$ws.Range("M68").Value = -907
$ws.Range("N68").Value = -3673

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# Row 71
$ws.Range("H71").Value = 1886.6666
$ws.Range("I71").Value = 1656
$ws.Range("J71").Value = 2175
$ws.Range("K71").Value = 8280
$ws.Range("L71").Value = 10875
$ws.Range("M71").Value = -4536
$ws.Range("N71").Value = -18363

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 10144.4
$ws.Range("I62").Value = 8840.799999999999
$ws.Range("J62").Value = 11448
$ws.Range("K62").Value = 8840.799999999999
$ws.Range("L62").Value = 11448
$ws.Range("M62").Value = -8216.799999999999
$ws.Range("N62").Value = -12696

# Row 65
$ws.Range("H65").Value = 10144.4
$ws.Range("I65").Value = 8840.799999999999
$ws.Range("J65").Value = 11448
$ws.Range("K65").Value = 44204
$ws.Range("L65").Value = 57240
$ws.Range("M65").Value = -41084
$ws.Range("N65").Value = -63480

# Row 81
$ws.Range("H81").Value = 8701.643
$ws.Range("I81").Value = 17636.834
$ws.Range("J81").Value = 2000.25
$ws.Range("K81").Value = 35273.668
$ws.Range("L81").Value = 4000.5
$ws.Range("M81").Value = -34212.668
$ws.Range("N81").Value = -6122.5

# Row 84
$ws.Range("H84").Value = 8701.643
$ws.Range("I84").Value = 17636.834
$ws.Range("J84").Value = 2000.25
$ws.Range("K84").Value = 176368.34
$ws.Range("L84").Value = 20002.5
$ws.Range("M84").Value = -171064.34
$ws.Range("N84").Value = -30610.5
